$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.435.50"
$ws.Range("E2").Value = "  +11.92%  "
$ws.Range("D3").Value = "1.829.13"
$ws.Range("E3").Value = "  +7.94%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'231.32"
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").Value = "'0.546"
$ws.Range("E6").Value = "  +4.35%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").Value = "'45.77"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("D10").Value = "'0.283"
$ws.Range("E10").Value = "  +5.87%  "
$ws.Range("E11").Value = "  +8.18%  "
$ws.Range("D12").Value = "'0.0932"
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("D13").Value = "2.094.12"
$ws.Range("E13").Value = "  +7.90%  "
$ws.Range("D14").Value = "1.829.74"
$ws.Range("E14").Value = "  +7.79%  "
$ws.Range("D15").Value = "'0.648"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("D16").Value = "34.417.26"
$ws.Range("E16").Value = "  +11.73%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'4.37"
$ws.Range("E17").Value = "  +8.54%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'10.28"
$ws.Range("E18").Value = "  -4.01%  "
$ws.Range("D19").Value = "'70.07"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").Value = "'260.26"
$ws.Range("E20").Value = "  +4.09%  "
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'10.56"
$ws.Range("E23").Value = "  +3.04%  "
$ws.Range("D24").Value = "'4.40"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'161.39"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").Value = "'16.86"
$ws.Range("E27").Value = "  +5.50%  "
$ws.Range("D28").Value = "'7.21"
$ws.Range("E28").Value = "  +6.44%  "
$ws.Range("E29").Value = "  +4.63%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'3.87"
$ws.Range("E31").Value = "  +10.75%  "
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("E33").Value = "  +6.99%  "
$ws.Range("D34").Value = "'3.58"
$ws.Range("E34").Value = "  +8.13%  "
$ws.Range("D35").Value = "1.580.01"
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("E36").Value = "  +5.74%  "
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0190"
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.633"
$ws.Range("E39").Value = "  +7.85%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'85.30"
$ws.Range("E40").Value = "  +5.67%  "
$ws.Range("D41").Value = "'2.87"
$ws.Range("E41").Value = "  +5.82%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").Value = "'0.921"
$ws.Range("E43").Value = "  +7.48%  "
$ws.Range("E44").Value = "  +5.69%  "
$ws.Range("D45").Value = "'0.0520"
$ws.Range("E45").Value = "  +3.32%  "
$ws.Range("D46").Value = "'1.06"
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("D47").Value = "1.984.76"
$ws.Range("E47").Value = "  +8.18%  "
$ws.Range("E48").Value = "  +5.90%  "
$ws.Range("D49").Value = "'53.07"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  +9.20%  "
